$d = $word.ActiveDocument

# Replace the placeholder text in the first paragraph, removing the
# trailing-space run by including it in the search text (Word merges the
# runs into a single run that keeps the first run's formatting).
$d.Content.Find.Execute("**ID__AFFARS_5325_topic_5__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5325_202__ID**", 2)

# Update the first paragraph's formatting: add a paragraph border (5pt
# space on each side) and widen the left indent from 120 to 225 twips
# (i.e. 6pt to 11.25pt).
$p = $d.Paragraphs(1)
$pf = $p.Range.ParagraphFormat
$pf.LeftIndent = 11.25

$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

Write-Output "done"
